# Fix for SLR update bug: a single reference with an oversized field (authors/
# abstract) used to abort the whole update. This edit represents the corrected
# re-run of the updater: (1) it adds a new "Misc. Data" column (J) to the
# References sheet, populated as empty for the existing rows, and (2) it
# rewrites the "Authors" column (E) for every existing reference row with the
# freshly (re-)fetched field values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column: "Misc. Data" (J) ---------------------------------------
$ws.Range("J1").Value = "Misc. Data"

# Existing rows get a blank Misc. Data value (column did not exist before).
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 10).Value = ""
}

# --- Re-write the Authors column (E) for each reference row -------------

$ws.Range("E2").Value = '[Daniel G.%Federman%NULL%1,                        Shaili%Gupta%NULL%2,                        Shaili%Gupta%NULL%0,                        Gary%Stack%NULL%1,                        Sheldon M.%Campbell%NULL%1,                        David R.%Peaper%NULL%2,                        David R.%Peaper%NULL%0,                        Louise M.%Dembry%NULL%1,                        Ann%Fisher%NULL%1,                        Asim F.%Tarabar%NULL%1,                        Michael%Kozal%NULL%1,                        Christopher B.%Ruser%NULL%1,                        Selvakumar%Subbian%NULL%8,                        Selvakumar%Subbian%NULL%0,                        Selvakumar%Subbian%NULL%0,                        Selvakumar%Subbian%NULL%0,                        Selvakumar%Subbian%NULL%0,                        Selvakumar%Subbian%NULL%0,                        Selvakumar%Subbian%NULL%0,                        Selvakumar%Subbian%NULL%0]'
$ws.Range("E3").Value = '[ B. M.%Berenger%null%1,                         K.% Fonseca%null%1,                         A. R.% Schneider%null%1,                         J.% Hu%null%1,                         N. % Zelyas%null%1,       B. M.%Berenger%null%1,       K.% Fonseca%null%1,       A. R.% Schneider%null%1,       J.% Hu%null%2,       N. % Zelyas%null%1]'
$ws.Range("E4").Value = '[Hélène%Péré%NULL%1,                        Isabelle%Podglajen%NULL%1,                        Maxime%Wack%NULL%1,                        Edouard%Flamarion%NULL%1,                        Tristan%Mirault%NULL%1,                        Guillaume%Goudot%NULL%1,                        Caroline%Hauw-Berlemont%NULL%2,                        Laetitia%Le%NULL%1,                        Eric%Caudron%NULL%1,                        Sophie%Carrabin%NULL%1,                        Julien%Rodary%NULL%1,                        Tatiana%Ribeyre%NULL%1,                        Laurent%Bélec%NULL%1,                        David%Veyer%NULL%1,                        Alexander J.%McAdam%NULL%7,                        Alexander J.%McAdam%NULL%0]'
$ws.Range("E5").Value = '[Yuan-Po%Tu%NULL%1,                        Rachel%Jennings%NULL%2,                        Rachel%Jennings%NULL%0,                        Brian%Hart%NULL%1,                        Gerard A.%Cangelosi%NULL%2,                        Gerard A.%Cangelosi%NULL%0,                        Rachel C.%Wood%NULL%1,                        Kevin%Wehber%NULL%2,                        Kevin%Wehber%NULL%0,                        Prateek%Verma%NULL%1,                        Deneen%Vojta%NULL%1,                        Ethan M.%Berke%NULL%1]'
$ws.Range("E6").Value = '[Cody%Callahan%NULL%1,                        Rose A.%Lee%NULL%1,                        Ghee Rye%Lee%NULL%1,                        Kate%Zulauf%NULL%1,                        James E.%Kirby%NULL%1,                        Ramy%Arnaout%NULL%1,                         C.%Callahan%null%1,                         R.% Lee%null%1,                         G.% Lee%null%1,                         K. E.% Zulauf%null%1,                         J. E.% Kirby%null%1,                         R. % Arnaout%null%1,                        C.%Callahan%null%1,                        R.% Lee%null%1,                        G.% Lee%null%1,                        K. E.% Zulauf%null%1,                        J. E.% Kirby%null%1,                        R. % Arnaout%null%1,       C.%Callahan%null%1,       R.% Lee%null%1,       G.% Lee%null%1,       K. E.% Zulauf%null%1,       J. E.% Kirby%null%1,       R. % Arnaout%null%1]'
$ws.Range("E7").Value = '[Sara B.%Griesemer%NULL%1,                        Greta%Van Slyke%NULL%1,                        Dylan%Ehrbar%NULL%1,                        Klemen%Strle%NULL%1,                        Tugba%Yildirim%NULL%1,                        Dominick A.%Centurioni%NULL%1,                        Anne C.%Walsh%NULL%1,                        Andrew K.%Chang%NULL%1,                        Michael J.%Waxman%NULL%2,                        Michael J.%Waxman%NULL%0,                        Kirsten%St. George%NULL%1,                        Alexander J.%McAdam%NULL%0,                        Alexander J.%McAdam%NULL%0,                        Alexander J.%McAdam%NULL%0]'
$ws.Range("E8").Value = '[N%Kojima%nkojima@ucla.edu%1,                        F%Turner%NULL%1,                        V%Slepnev%NULL%1,                        A%Bacelar%NULL%1,                        L%Deming%NULL%1,                        S%Kodeboyina%NULL%1,                        J D%Klausner%NULL%1]'
$ws.Range("E9").Value = '[Denise J.%McCulloch%NULL%1,                        Ashley E.%Kim%NULL%1,                        Naomi C.%Wilcox%NULL%1,                        Jennifer K.%Logue%NULL%1,                        Alex L.%Greninger%NULL%1,                        Janet A.%Englund%NULL%1,                        Helen Y.%Chu%NULL%0]'
$ws.Range("E10").Value = '[K. E.%Hanson%NULL%1,                        A. P.%Barker%NULL%2,                        A. P.%Barker%NULL%0,                        D. R.%Hillyard%NULL%1,                        N.%Gilmore%NULL%1,                        J. W.%Barrett%NULL%1,                        R. R.%Orlandi%NULL%1,                        S. M.%Shakir%NULL%1,                        Alexander J.%McAdam%NULL%0,                        Alexander J.%McAdam%NULL%0]'
$ws.Range("E11").Value = '[Swetha%Pinninti%spinninti@peds.uab.edu%1,                        Connie%Trieu%NULL%1,                        Sunil K%Pati%NULL%1,                        Misty%Latting%NULL%1,                        Joshua%Cooper%NULL%1,                        Maria C%Seleme%NULL%1,                        Sushma%Boppana%NULL%1,                        Nitin%Arora%NULL%1,                        William J%Britt%NULL%1,                        Suresh B%Boppana%NULL%1]'
$ws.Range("E12").Value = '[Min%Liu%NULL%1,                        Qianyuan%Li%NULL%1,                        Jun%Zhou%NULL%1,                        Wen%Ai%NULL%1,                        Xiaoling%Zheng%NULL%1,                        Jingjing%Zeng%NULL%1,                        Yuwen%Liu%NULL%1,                        Xiying%Xiang%NULL%1,                        Rong%Guo%NULL%1,                        Xiaoyin%Li%NULL%1,                        Xiandi%Wu%NULL%1,                        Haiying%Xu%NULL%1,                        Ling%Jiang%NULL%1,                        Huaqin%Zhang%NULL%1,                        Jing%Chen%NULL%2,                        Lili%Tian%NULL%1,                        Jun%Luo%NULL%1,                        Chunhua%Luo%NULL%1]'
